$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("W2").Value = 10
$ws.Range("AC2").Value = 8
$ws.Range("BC2").Value = 126
$ws.Range("BD2").Value = 126
$ws.Range("O3").Value = 1.67
$ws.Range("P3").Value = 2.1
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("R7").Value = 1.62
$ws.Range("R8").Value = 1.57
$ws.Range("G10").Value = 3
$ws.Range("I10").Value = 2.6
$ws.Range("J10").Value = 3.75
$ws.Range("L10").Value = 3.4
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("X10").Value = 13
$ws.Range("Y10").Value = 12
$ws.Range("Z10").Value = 34
$ws.Range("AD10").Value = 5.5
$ws.Range("AH10").Value = 6.5
$ws.Range("AI10").Value = 11
$ws.Range("AK10").Value = 26
$ws.Range("AL10").Value = 26
$ws.Range("AN10").Value = 4.75
$ws.Range("AO10").Value = 19
$ws.Range("AQ10").Value = 67
$ws.Range("AW10").Value = 4.33
$ws.Range("AX10").Value = 15
$ws.Range("AY10").Value = 29
$ws.Range("BA10").Value = 81
$ws.Range("AC11").Value = 9
$ws.Range("AD11").Value = 6.5
$ws.Range("G12").Value = 2.3
$ws.Range("H12").Value = 2.9
$ws.Range("I12").Value = 3.4
$ws.Range("J12").Value = 3.2
$ws.Range("L12").Value = 4
$ws.Range("S12").Value = 1.57
$ws.Range("T12").Value = 2.25
$ws.Range("W12").Value = 6.5
$ws.Range("X12").Value = 10
$ws.Range("Z12").Value = 23
$ws.Range("AA12").Value = 23
$ws.Range("AH12").Value = 8
$ws.Range("AI12").Value = 15
$ws.Range("AK12").Value = 34
$ws.Range("AN12").Value = 4.33
$ws.Range("AO12").Value = 15
$ws.Range("AQ12").Value = 51
$ws.Range("AT12").Value = 2.25
$ws.Range("AX12").Value = 19
$ws.Range("M13").Value = 1.11
$ws.Range("N13").Value = 6.5
$ws.Range("AH13").Value = 9.5
$ws.Range("AI13").Value = 21
$ws.Range("AJ13").Value = 17
$ws.Range("AK13").Value = 51
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 15
$ws.Range("O14").Value = 1.2
$ws.Range("P14").Value = 4.33
$ws.Range("Q14").Value = 1.67
$ws.Range("R14").Value = 2.15
$ws.Range("BD14").Value = 151
$ws.Range("G17").Value = 3.1
$ws.Range("I17").Value = 2.4
$ws.Range("L17").Value = 3.25
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5
$ws.Range("W17").Value = 7.5
$ws.Range("X17").Value = 13
$ws.Range("Y17").Value = 12
$ws.Range("AA17").Value = 29
$ws.Range("AH17").Value = 6.5
$ws.Range("AK17").Value = 23
$ws.Range("AO17").Value = 19
$ws.Range("G18").Value = 3.4
$ws.Range("J18").Value = 4.2
$ws.Range("N18").Value = 7
$ws.Range("O18").Value = 1.47
$ws.Range("P18").Value = 2.8
$ws.Range("Q18").Value = 2.37
$ws.Range("R18").Value = 1.62
$ws.Range("S18").Value = 1.55
$ws.Range("T18").Value = 2.45
$ws.Range("U18").Value = 1.98
$ws.Range("V18").Value = 1.78
$ws.Range("W18").Value = 8.5
$ws.Range("X18").Value = 18.5
$ws.Range("AA18").Value = 40
$ws.Range("AB18").Value = 55
$ws.Range("AC18").Value = 7
$ws.Range("AE18").Value = 18
$ws.Range("AF18").Value = 110
$ws.Range("AG18").Value = 1250
$ws.Range("AH18").Value = 6.2
$ws.Range("AJ18").Value = 9.75
$ws.Range("AL18").Value = 22
$ws.Range("AM18").Value = 40
$ws.Range("AO18").Value = 21
$ws.Range("AP18").Value = 32
$ws.Range("AQ18").Value = 120
$ws.Range("AR18").Value = 175
$ws.Range("AS18").Value = 500
$ws.Range("AT18").Value = 2.42
$ws.Range("AU18").Value = 7.9
$ws.Range("AV18").Value = 90
$ws.Range("AW18").Value = 3.9
$ws.Range("BA18").Value = 110
$ws.Range("G19").Value = 2.1
$ws.Range("H19").Value = 3.1
$ws.Range("I19").Value = 3.75
$ws.Range("J19").Value = 2.75
$ws.Range("L19").Value = 4.33
$ws.Range("M19").Value = 1.08
$ws.Range("N19").Value = 8
$ws.Range("U19").Value = 1.83
$ws.Range("V19").Value = 1.83
$ws.Range("W19").Value = 7
$ws.Range("X19").Value = 9.5
$ws.Range("Z19").Value = 19
$ws.Range("AA19").Value = 19
$ws.Range("AG19").Value = 301
$ws.Range("AJ19").Value = 13
$ws.Range("AO19").Value = 12
$ws.Range("AR19").Value = 67
$ws.Range("AY19").Value = 29
$ws.Range("AZ19").Value = 67

Write-Output "Applied 131 cell updates"
